$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 70: the date must be written as literal text ("2025/10/06"),
# matching the existing rows, not auto-converted to a date serial number.
# Briefly force Text format so Excel's auto-detect leaves the string alone,
# then clear the format again so the new cell has no style override
# (matching the rest of the sheet, which carries no per-cell style either).
$ws.Range("A70").NumberFormat = "@"
$ws.Range("A70").Value = "2025/10/06"
$ws.Range("A70").ClearFormats()

$ws.Range("B70").Value = "月"
$ws.Range("C70").Value = 21
$ws.Range("D70").Value = 6
